$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage (no auto-number-coercion) for the edited range, without altering cell style.
$editRange = $ws.Range("B2:E51")
$editRange.NumberFormat = "@"

$ws.Range("D2").Value = '55.324.45'
$ws.Range("E2").Value = '  -1.85%  '
$ws.Range("D3").Value = '2.350.33'
$ws.Range("E3").Value = '  -5.52%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").Value = '476.18'
$ws.Range("E5").Value = '  -2.86%  '
$ws.Range("D6").Value = '146.16'
$ws.Range("E6").Value = '  -1.49%  '
$ws.Range("E7").Value = '  +23.08%  '
$ws.Range("D8").Value = '0.999'
$ws.Range("E8").Value = '  +0.22%  '
$ws.Range("D9").Value = '2.355.19'
$ws.Range("E9").Value = '  -5.73%  '
$ws.Range("D10").Value = '0.0964'
$ws.Range("E10").Value = '  -1.40%  '
$ws.Range("D11").Value = '5.47'
$ws.Range("E11").Value = '  -6.18%  '
$ws.Range("E12").Value = '  -2.10%  '
$ws.Range("D13").Value = '0.124'
$ws.Range("E13").Value = '  +0.61%  '
$ws.Range("D14").Value = '2.759.87'
$ws.Range("E14").Value = '  -5.39%  '
$ws.Range("D15").Value = '55.237.40'
$ws.Range("E15").Value = '  -1.88%  '
$ws.Range("D16").Value = '20.04'
$ws.Range("E16").Value = '  -5.67%  '
$ws.Range("D17").Value = '0.0000129'
$ws.Range("E17").Value = '  -5.54%  '
$ws.Range("D18").Value = '2.360.43'
$ws.Range("E18").Value = '  -5.12%  '
$ws.Range("E19").Value = '  +0.17%  '
$ws.Range("D20").Value = '314.51'
$ws.Range("E20").Value = '  -1.90%  '
$ws.Range("D21").Value = '9.57'
$ws.Range("E21").Value = '  -5.68%  '
$ws.Range("E22").Value = '  -0.07%  '
$ws.Range("E23").Value = '  -3.01%  '
$ws.Range("D24").Value = '56.65'
$ws.Range("E24").Value = '  -3.63%  '
$ws.Range("D25").Value = '0.999'
$ws.Range("E25").Value = '  +0.00%  '
$ws.Range("D26").Value = '0.394'
$ws.Range("E26").Value = '  -4.80%  '
$ws.Range("D27").Value = '0.151'
$ws.Range("E27").Value = '  -7.33%  '
$ws.Range("D28").Value = '2.454.38'
$ws.Range("E28").Value = '  -5.08%  '
$ws.Range("D29").Value = '7.09'
$ws.Range("E29").Value = '  -7.66%  '
$ws.Range("E30").Value = '  +0.16%  '
$ws.Range("D31").Value = '0.0₃0744'
$ws.Range("E31").Value = '  -6.73%  '
$ws.Range("B32").Value = 'Monero'
$ws.Range("C32").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D32").Value = '145.96'
$ws.Range("E32").Value = '  -2.34%  '
$ws.Range("B33").Value = 'EthereumClassic'
$ws.Range("C33").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D33").Value = '18.08'
$ws.Range("E33").Value = '  -1.26%  '
$ws.Range("E34").Value = '  -2.62%  '
$ws.Range("D35").Value = '5.09'
$ws.Range("E35").Value = '  -2.62%  '
$ws.Range("E36").Value = '  -4.43%  '
$ws.Range("D37").Value = '3.59'
$ws.Range("E37").Value = '  -4.60%  '
$ws.Range("D38").Value = '0.812'
$ws.Range("E38").Value = '  -6.16%  '
$ws.Range("D39").Value = '0.102'
$ws.Range("E39").Value = '  +10.33%  '
$ws.Range("E40").Value = '  -1.16%  '
$ws.Range("D41").Value = '0.998'
$ws.Range("E41").Value = '  +0.32%  '
$ws.Range("D42").Value = '1.33'
$ws.Range("E42").Value = '  -0.77%  '
$ws.Range("D43").Value = '3.39'
$ws.Range("E43").Value = '  -4.30%  '
$ws.Range("E44").Value = '  -5.19%  '
$ws.Range("D45").Value = '0.0518'
$ws.Range("E45").Value = '  -7.06%  '
$ws.Range("E46").Value = '  -0.36%  '
$ws.Range("D47").Value = '251.53'
$ws.Range("E47").Value = '  -4.36%  '
$ws.Range("D48").Value = '0.0220'
$ws.Range("E48").Value = '  -3.77%  '
$ws.Range("E49").Value = '  -7.74%  '
$ws.Range("B50").Value = 'Maker'
$ws.Range("C50").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D50").Value = '1.798.51'
$ws.Range("E50").Value = '  -4.70%  '
$ws.Range("B51").Value = 'EnergySwap'
$ws.Range("C51").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D51").Value = '16.69'
$ws.Range("E51").Value = '  -5.11%  '

# Restore the default (Normal) style so cell formatting matches the original workbook.
$editRange.Style = "Normal"
